# update timesleep to not take too long
# Trims each sheet's product list from 5 products down to 3 products,
# replacing the first three rows' data with new product rows and deleting
# the trailing two rows (shifting the Total row up on the Summary sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Summary"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop the last two product rows (old rows 5 & 6); the Total row (old
# row 7) shifts up to row 5 and its SUM formulas auto-adjust.
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(5).Delete()

# Overwrite the remaining three product rows with the new data.
$ws1.Range("A2").Value = "DRÖNA / KAVALKAD"
$ws1.Range("B2").Value = "402.493.53"
$ws1.Range("C2").Value = "33x38x33 cm"
$ws1.Range("D2").Value = 129
$ws1.Range("E2").Value = 95.50533003122644

$ws1.Range("A3").Value = "VARDAGEN / SANDVIVA"
$ws1.Range("B3").Value = "002.947.24"
$ws1.Range("C3").Value = "23 cm"
$ws1.Range("D3").Value = 449
$ws1.Range("E3").Value = 95.50533003122644

$ws1.Range("A4").Value = "SITTBRUNN / MÅLA"
$ws1.Range("B4").Value = "805.394.83"
$ws1.Range("C4").Value = "1 m"
$ws1.Range("D4").Value = 49
$ws1.Range("E4").Value = 25.13298158716485

# Column A narrows from 25 to 21 characters.
$ws1.Columns.Item(1).ColumnWidth = 21 - 5/6

# ---------------------------------------------------------------------
# Sheet 2: "Czech Data"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(5).Delete()
$ws2.Rows.Item(5).Delete()

$ws2.Range("A2").Value = "DRÖNA"
$ws2.Range("B2").Value = 25.66348913928281
$ws2.Range("C2").Value = 129
$ws2.Range("D2").Value = "402.493.53"
$ws2.Range("E2").Value = "Krabice, červená,"
$ws2.Range("F2").Value = "33x38x33 cm"

$ws2.Range("A3").Value = "VARDAGEN"
$ws2.Range("B3").Value = 89.32485754680606
$ws2.Range("C3").Value = 449
$ws2.Range("D3").Value = "002.947.24"
$ws2.Range("E3").Value = "Nůž na chléb, tmavě šedá,"
$ws2.Range("F3").Value = "23 cm"

$ws2.Range("A4").Value = "SITTBRUNN"
$ws2.Range("B4").Value = 9.748147037401997
$ws2.Range("C4").Value = 49
$ws2.Range("D4").Value = "805.394.83"
$ws2.Range("E4").Value = "USB-A na USB-C, světle žlutá,"
$ws2.Range("F4").Value = "1 m"

# Column A narrows from 16 to 14, column E from 50 to 31.
$ws2.Columns.Item(1).ColumnWidth = 14 - 5/6
$ws2.Columns.Item(5).ColumnWidth = 31 - 5/6

# ---------------------------------------------------------------------
# Sheet 3: "Poland Data"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(5).Delete()
$ws3.Rows.Item(5).Delete()

$ws3.Range("A2").Value = "KAVALKAD"
$ws3.Range("B2").Value = 19
$ws3.Range("C2").Value = 95.50533003122644
$ws3.Range("D2").Value = "002.677.06"
$ws3.Range("E2").Value = "Patelnia, czarny,"
$ws3.Range("F2").Value = "24 cm"

$ws3.Range("A3").Value = "SANDVIVA"
$ws3.Range("B3").Value = 19
$ws3.Range("C3").Value = 95.50533003122644
$ws3.Range("D3").Value = "104.643.82"
$ws3.Range("E3").Value = "Rękawica, silikon/niebieski"
# F3 stays "Not available" (unchanged by the diff).

$ws3.Range("A4").Value = "MÅLA"
$ws3.Range("B4").Value = 5
$ws3.Range("C4").Value = 25.13298158716485
$ws3.Range("D4").Value = "904.565.90"
$ws3.Range("E4").Value = "Nożyczki"
$ws3.Range("F4").Value = "Not available"

# Column A narrows from 16 to 14, column D from 15 to 14, column E from 44 to 29.
$ws3.Columns.Item(1).ColumnWidth = 14 - 5/6
$ws3.Columns.Item(4).ColumnWidth = 14 - 5/6
$ws3.Columns.Item(5).ColumnWidth = 29 - 5/6
